# ---- Part 1: Shared-string text edits (Volume number + report date range) ----
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A8: "Volume 32   Number  23" -> "...24" (last run "23" -> "24")
$rA8 = $ws.Range("A8")
$cA8 = $rA8.Characters(21, 2)
$cA8.Text = "24"

# C9: "Report Covering the Week  6/2/2025  Through  6/8/2025"
#     -> "...6/9/2025  Through  6/15/2025"
$rC9 = $ws.Range("C9")
$c1C9 = $rC9.Characters(27, 8)
$c1C9.Text = "6/9/2025"
$c2C9 = $rC9.Characters(46, 8)
$c2C9.Text = "6/15/2025"

# ---- Part 2: Plain numeric value updates (style/type unchanged) ----
$ws.Range("M14").Value = -50

$ws.Range("G15").Value = 7
$ws.Range("H15").Value = -14.285714285714
$ws.Range("I15").Value = 35
$ws.Range("J15").Value = 31
$ws.Range("K15").Value = 12.903225806451
$ws.Range("L15").Value = 75
$ws.Range("M15").Value = 52.173913043478
$ws.Range("N15").Value = -16.666666666666

$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -60
$ws.Range("F16").Value = 34
$ws.Range("G16").Value = 21
$ws.Range("H16").Value = 61.904761904761
$ws.Range("I16").Value = 136
$ws.Range("J16").Value = 144
$ws.Range("K16").Value = -5.555555555555
$ws.Range("L16").Value = -4.225352112676
$ws.Range("M16").Value = -22.727272727272
$ws.Range("N16").Value = -77.371048252911

$ws.Range("C17").Value = 18
$ws.Range("D17").Value = 23
$ws.Range("E17").Value = -21.739130434782
$ws.Range("F17").Value = 89
$ws.Range("G17").Value = 77
$ws.Range("H17").Value = 15.584415584415
$ws.Range("I17").Value = 457
$ws.Range("J17").Value = 450
$ws.Range("K17").Value = 1.555555555555
$ws.Range("L17").Value = 9.069212410501
$ws.Range("M17").Value = 97.835497835497
$ws.Range("N17").Value = -10.03937007874

$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 17
$ws.Range("G18").Value = 32
$ws.Range("H18").Value = -46.875
$ws.Range("I18").Value = 157
$ws.Range("J18").Value = 140
$ws.Range("K18").Value = 12.142857142857
$ws.Range("L18").Value = -1.25786163522
$ws.Range("M18").Value = -35.918367346938
$ws.Range("N18").Value = -89.718402095612

$ws.Range("C19").Value = 22
$ws.Range("D19").Value = 49
$ws.Range("E19").Value = -55.102040816326
$ws.Range("F19").Value = 105
$ws.Range("G19").Value = 139
$ws.Range("H19").Value = -24.460431654676
$ws.Range("I19").Value = 596
$ws.Range("J19").Value = 745
$ws.Range("K19").Value = -20
$ws.Range("L19").Value = -16.875871687587
$ws.Range("M19").Value = 47.524752475247
$ws.Range("N19").Value = -15.700141442715

$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 8
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 16
$ws.Range("G20").Value = 39
$ws.Range("H20").Value = -58.974358974359
$ws.Range("I20").Value = 84
$ws.Range("J20").Value = 137
$ws.Range("K20").Value = -38.686131386861
$ws.Range("L20").Value = -55.555555555555
$ws.Range("M20").Value = -36.363636363636
$ws.Range("N20").Value = -96.159122085048

$ws.Range("C21").Value = 49
$ws.Range("D21").Value = 93
$ws.Range("E21").Value = -47.311827956989
$ws.Range("F21").Value = 267
$ws.Range("G21").Value = 315
$ws.Range("H21").Value = -15.238095238095
$ws.Range("I21").Value = 1467
$ws.Range("J21").Value = 1648
$ws.Range("K21").Value = -10.983009708737
$ws.Range("L21").Value = -11.41304347826
$ws.Range("M21").Value = 20.74074074074
$ws.Range("N21").Value = -73.723804406233

$ws.Range("C23").Value = 2
$ws.Range("F23").Value = 10
$ws.Range("H23").Value = 25
$ws.Range("I23").Value = 68
$ws.Range("K23").Value = 38.775510204081
$ws.Range("L23").Value = -4.225352112676
$ws.Range("M23").Value = 172

$ws.Range("C24").Value = 73
$ws.Range("D24").Value = 85
$ws.Range("E24").Value = -14.117647058823
$ws.Range("F24").Value = 305
$ws.Range("G24").Value = 291
$ws.Range("H24").Value = 4.810996563573
$ws.Range("I24").Value = 1951
$ws.Range("J24").Value = 1950
$ws.Range("K24").Value = 0.051282051282
$ws.Range("L24").Value = 3.172924378635
$ws.Range("M24").Value = 22.627278441231

$ws.Range("C25").Value = 34
$ws.Range("D25").Value = 42
$ws.Range("E25").Value = -19.047619047619
$ws.Range("F25").Value = 189
$ws.Range("G25").Value = 153
$ws.Range("H25").Value = 23.529411764705
$ws.Range("I25").Value = 1166
$ws.Range("J25").Value = 1071
$ws.Range("K25").Value = 8.870214752567
$ws.Range("L25").Value = 27.571115973741

$ws.Range("C26").Value = 47
$ws.Range("D26").Value = 43
$ws.Range("E26").Value = 9.302325581395
$ws.Range("F26").Value = 155
$ws.Range("G26").Value = 162
$ws.Range("H26").Value = -4.32098765432
$ws.Range("I26").Value = 866
$ws.Range("J26").Value = 839
$ws.Range("K26").Value = 3.218116805721
$ws.Range("L26").Value = 8.793969849246
$ws.Range("M26").Value = -4.203539823008

$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 8
$ws.Range("G27").Value = 9
$ws.Range("H27").Value = -11.111111111111
$ws.Range("I27").Value = 45
$ws.Range("J27").Value = 54
$ws.Range("K27").Value = -16.666666666666
$ws.Range("L27").Value = 50

$ws.Range("F28").Value = 8
$ws.Range("G28").Value = 16
$ws.Range("H28").Value = -50
$ws.Range("I28").Value = 98
$ws.Range("J28").Value = 87
$ws.Range("K28").Value = 12.643678160919
$ws.Range("L28").Value = -3.92156862745

$ws.Range("J29").Value = 7
$ws.Range("K29").Value = -57.142857142857
$ws.Range("M29").Value = -81.25
$ws.Range("N29").Value = -94.117647058823

$ws.Range("J30").Value = 7
$ws.Range("K30").Value = -57.142857142857
$ws.Range("M30").Value = -78.571428571428
$ws.Range("N30").Value = -92.682926829268

$ws.Range("D31").Value = 2
$ws.Range("G31").Value = 4
$ws.Range("H31").Value = -75
$ws.Range("I31").Value = 5
$ws.Range("J31").Value = 13
$ws.Range("K31").Value = -61.538461538461
$ws.Range("L31").Value = -28.571428571428

$ws.Range("H33").Value = -100

# ---- Part 3: Cells switching from NUMBER style to TEXT/placeholder style ("0" / "***.*") ----
# Use a leading apostrophe to force text interpretation, then copy the General/text
# format (style index 13 in the source file) from a known placeholder cell (C14).
$ws.Range("D23").Value = "'0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4122) | Out-Null
$ws.Range("E23").Value = "'***.*"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("E23").PasteSpecial(-4122) | Out-Null

$ws.Range("D33").Value = "'0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D33").PasteSpecial(-4122) | Out-Null
$ws.Range("E33").Value = "'***.*"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("E33").PasteSpecial(-4122) | Out-Null
$ws.Range("F33").Value = "'0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("F33").PasteSpecial(-4122) | Out-Null

# ---- Part 4: Cells switching from TEXT/placeholder style to NUMBER style ----
# Write the numeric value, then copy the correct numeric format (style 14 = integer,
# style 15 = one-decimal percent) from known source cells (C39 / K39).
$ws.Range("D29").Value = 1
$ws.Range("C39").Copy() | Out-Null
$ws.Range("D29").PasteSpecial(-4122) | Out-Null
$ws.Range("E29").Value = -100
$ws.Range("K39").Copy() | Out-Null
$ws.Range("E29").PasteSpecial(-4122) | Out-Null
$ws.Range("G29").Value = 1
$ws.Range("C39").Copy() | Out-Null
$ws.Range("G29").PasteSpecial(-4122) | Out-Null
$ws.Range("H29").Value = -100
$ws.Range("K39").Copy() | Out-Null
$ws.Range("H29").PasteSpecial(-4122) | Out-Null

$ws.Range("D30").Value = 1
$ws.Range("C39").Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4122) | Out-Null
$ws.Range("E30").Value = -100
$ws.Range("K39").Copy() | Out-Null
$ws.Range("E30").PasteSpecial(-4122) | Out-Null
$ws.Range("G30").Value = 1
$ws.Range("C39").Copy() | Out-Null
$ws.Range("G30").PasteSpecial(-4122) | Out-Null
$ws.Range("H30").Value = -100
$ws.Range("K39").Copy() | Out-Null
$ws.Range("H30").PasteSpecial(-4122) | Out-Null

$ws.Range("F31").Value = 1
$ws.Range("C39").Copy() | Out-Null
$ws.Range("F31").PasteSpecial(-4122) | Out-Null
